$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-11 (years 2000-2009), shifting rows 12-22 (years 2010-2020) up to rows 2-12
$ws.Range("A2:E11").EntireRow.Delete()
